$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.203.66"
$ws.Range("E2").Value = "  -0.41%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.536.11"
$ws.Range("E3").Value = "  +3.21%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.74"
$ws.Range("E5").Value = "  +1.67%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.42"
$ws.Range("E6").Value = "  +0.52%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.535.64"
$ws.Range("E7").Value = "  +3.18%  "
$ws.Range("E9").Value = "  -1.05%  "
$ws.Range("E10").Value = "  +2.96%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.82"
$ws.Range("E11").Value = "  -5.96%  "
$ws.Range("E12").Value = "  +3.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.138.23"
$ws.Range("E13").Value = "  +3.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000186"
$ws.Range("E14").Value = "  +2.85%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.28"
$ws.Range("E15").Value = "  +3.88%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.536.58"
$ws.Range("E16").Value = "  +2.98%  "
$ws.Range("E17").Value = "  +1.45%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "65.234.89"
$ws.Range("E18").Value = "  -0.26%  "
$ws.Range("E19").Value = "  +4.76%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.96"
$ws.Range("E20").Value = "  +1.72%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.36"
$ws.Range("E21").Value = "  +5.27%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "393.60"
$ws.Range("E22").Value = "  +0.67%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.575"
$ws.Range("E23").Value = "  +3.65%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.678.99"
$ws.Range("E24").Value = "  +3.36%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "73.88"
$ws.Range("E25").Value = "  +1.00%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("E27").Value = "  +8.42%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.80"
$ws.Range("E28").Value = "  +8.98%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.28"
$ws.Range("E30").Value = "  +2.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.17"
$ws.Range("E31").Value = "  -0.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.548.11"
$ws.Range("E32").Value = "  +3.52%  "
$ws.Range("E34").Value = "  +3.77%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.146"
$ws.Range("E35").Value = "  +0.87%  "
$ws.Range("E36").Value = "  +9.87%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.99"
$ws.Range("E37").Value = "  +2.00%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "169.18"
$ws.Range("E38").Value = "  -1.94%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.55"
$ws.Range("E39").Value = "  +5.68%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.00"
$ws.Range("E40").Value = "  +5.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0803"
$ws.Range("E41").Value = "  +5.26%  "
$ws.Range("E42").Value = "  +0.34%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.14"
$ws.Range("E43").Value = "  +16.41%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.78"
$ws.Range("E44").Value = "  -1.89%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.999"
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.43"
$ws.Range("E46").Value = "  +0.34%  "
$ws.Range("B47").Value = "ONDO"
$ws.Range("C47").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.19"
$ws.Range("E47").Value = "  +7.22%  "
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.68"
$ws.Range("E48").Value = "  +4.38%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.81"
$ws.Range("E49").Value = "  +4.36%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.383.12"
$ws.Range("E50").Value = "  +7.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "306.18"
$ws.Range("E51").Value = "  +7.13%  "
